# Fix wrong use of arrow
#
# Target: the "Elbow Connector 106" shape (Id=19) on slide 5 of the deck,
# the blue elbow connector whose tail previously ended in a connection to
# shape Id=18 ("Rectangle 62"/idx 1).
#
# Changes:
#   - detach the tail end from its connected shape (it was wrongly glued
#     to "Rectangle 62")
#   - flip the connector horizontally (it already had a vertical flip)
#   - change the tail (end) arrowhead from a plain "med/med" open arrow to
#     a bold "lg/lg" triangle arrowhead

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Locate the connector by its stable shape Id/Name (more robust than a
# hard-coded positional index into the Shapes collection).
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 19 -and $candidate.Name -eq "Elbow Connector 106") {
        $targetShape = $candidate
        break
    }
}

if ($targetShape -eq $null) {
    # Fall back to the known positional index if the Id-based lookup
    # somehow fails.
    $targetShape = $s.Shapes.Item(16)
}

# Detach the end point of the connector from "Rectangle 62" (shape Id 18).
try {
    $targetShape.ConnectorFormat.EndDisconnect()
} catch {
    # Older/limited hosts may not support live re-wiring of connectors;
    # ignore and continue with the formatting changes below.
}

# Add a horizontal flip (the shape already has a vertical flip applied).
$targetShape.HorizontalFlip = -1

# Swap the tail arrowhead from the default "arrow" (open, medium/medium)
# to a large triangle arrowhead.
$targetShape.Line.EndArrowheadStyle = 2   # msoArrowheadTriangle
$targetShape.Line.EndArrowheadWidth = 3   # msoArrowheadWide
$targetShape.Line.EndArrowheadLength = 3  # msoArrowheadLong
